$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fill in the new time-log entry in row 24 -----------------------------
# Set the interruption (D) first so the shared formula in E24 picks up the
# new value when it recalculates, then the rest of the row.
$ws.Range("D24").Value = 15
$ws.Range("A24").Value = 41885
$ws.Range("B24").Value = 0.92013888888888884
$ws.Range("C24").Value = 0.9868055555555556
$ws.Range("F24").Value = "Coding"
$ws.Calculate()

# --- Update the pie chart --------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

# Give the single series an explicit name (Sheet1!$I$1 = "Hours") so it shows
# up as the chart/series title.
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(Sheet1!`$I`$1,Sheet1!`$H`$2:`$H`$5,Sheet1!`$I`$2:`$I`$5,1)"

# Turn the chart title back on.
$chart.HasTitle = $true
$chart.ChartTitle.Text = "Hours"

# Show category names in the data labels.
$dl = $series.DataLabels()
$dl.ShowCategoryName = $true
$dl.ShowPercentage = $true
$dl.ShowValue = $false
$dl.ShowSeriesName = $false
$dl.ShowLegendKey = $false
$dl.ShowBubbleSize = $false

# Resize / reposition the chart.
$co.Left = 678.5693700787401
$co.Top = 15.37488188976378
$co.Width = 514.0
$co.Height = 253.87503937007875

# --- Restore the user's last selection -------------------------------------
$ws.Range("C25").Select()
